# Generate Report for Handback
# Mirrors the localization-status report refresh: the "Ready for handoff"
# status becomes "Handed back: in sync with en-US", handback timestamps are
# recorded for zh-cn/de-de, and the per-language sheets grow two extra
# columns worth of data (Latest Target File / Latest Handback File) plus a
# hyperlink to the generated target file.

$wb = $excel.ActiveWorkbook

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/923146179a1ed2a5d146ad92c359fdf5486a09bb/e2e/"
$targetFile = "3126d764-c01c-45fa-9a8e-a48247a84523.md"
$targetUrl = $ghBase + $targetFile

# ---------------------------------------------------------------------
# Overview sheet: both rows are now "Handed back", in sync with en-US.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------
# zh-cn sheet: fill in the generated target file + handback file, add the
# hyperlink to the target file, and stamp the handback datetime.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("J2").Value = "3126d764-c01c-45fa-9a8e-a48247a84523.680a4c83cbb864c2316e331da7c9f9ed0289fe78.zh-cn.xlf"
$zhcn.Range("J3").Value = "3126d764-c01c-45fa-9a8e-a48247a84523.680a4c83cbb864c2316e331da7c9f9ed0289fe78.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $targetUrl, "", "", $targetFile)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $targetUrl, "", "", $targetFile)

$zhcn.Range("K2").Value = "2016-08-26 11:00:46"
$zhcn.Range("K3").Value = "2016-08-26 11:00:46"

$zhcn.Columns.Item(3).ColumnWidth = 29.144371396019366
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: same shape of change, different handback file + datetime.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("J2").Value = "3126d764-c01c-45fa-9a8e-a48247a84523.680a4c83cbb864c2316e331da7c9f9ed0289fe78.de-de.xlf"
$dede.Range("J3").Value = "3126d764-c01c-45fa-9a8e-a48247a84523.680a4c83cbb864c2316e331da7c9f9ed0289fe78.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("I2"), $targetUrl, "", "", $targetFile)
$dede.Hyperlinks.Add($dede.Range("I3"), $targetUrl, "", "", $targetFile)

$dede.Range("K2").Value = "2016-08-26 11:00:53"
$dede.Range("K3").Value = "2016-08-26 11:00:53"

$dede.Columns.Item(3).ColumnWidth = 29.144371396019366
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
